$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "compatible Beta" product line used on the z0bug.invoice_ZI_6 / ZI_8
# invoices was renamed to "compatible Phi" (column E = name).
$ws.Range("E51").Value = "Prodotto compatibile Phi"
$ws.Range("E54").Value = "Prodotto compatibile Phi"

# Append a brand-new invoice (z0bug.invoice_ZI_10) made of two lines, built
# by duplicating the ZI_8 pair of rows (54-55: a "compatible product" line
# followed by its commercial-information line) down to rows 58-59 and then
# overwriting the values that differ.
$ws.Range("A54:K54").Copy($ws.Range("A58:K58"))
$ws.Range("M54:N54").Copy($ws.Range("M58:N58"))
$ws.Range("A55:K55").Copy($ws.Range("A59:K59"))

$ws.Range("A58").Value = "z0bug.invoice_ZI_10_1"
$ws.Range("B58").Value = "z0bug.invoice_ZI_10"
$ws.Range("D58").Value = "z0bug.product_product_6"
$ws.Range("E58").Value = "Prodotto Phi"
$ws.Range("F58").Value = 97
$ws.Range("G58").Value = "external.610100"
$ws.Range("H58").Value = 2.8081
$ws.Range("I58").Value = "z0bug.tax_a8aa"
$ws.Range("K58").Value = 1
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()

$ws.Range("A59").Value = "z0bug.invoice_ZI_10_2"
$ws.Range("B59").Value = "z0bug.invoice_ZI_10"
$ws.Range("D59").Value = "z0bug.product_product_7"
$ws.Range("E59").Value = "Prodotto Gamma"
$ws.Range("F59").Value = 91
$ws.Range("G59").Value = "external.610100"
$ws.Range("H59").Value = 0.3483
$ws.Range("I59").Value = "z0bug.tax_a8aa"
$ws.Range("K59").Value = 1

# Move the selection/view down to the newly added rows, like the author
# did after appending the new lines (header row stays frozen).
$ws.Range("A58").Select()
